$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$sectionName = "Banco de dados MySQL e Implantação no Heroku"
$lessonName  = "Instalando e preparando o Heroku"

# --- Extend the "Tabela1" table by 4 rows (130-133), matching how the
# author continued the table in Excel when annotating a new lesson ---
$tbl = $ws.ListObjects.Item(1)
$tbl.ListRows.Add() | Out-Null
$tbl.ListRows.Add() | Out-Null
$tbl.ListRows.Add() | Out-Null
$tbl.ListRows.Add() | Out-Null

# --- Set values in the same order the original author entered them, so that
# new shared-string table entries line up with the reference workbook ---

# 1) Lesson name column (E) first -> creates shared string "Instalando e preparando o Heroku"
$ws.Cells.Item(130, 5).Value = $lessonName
$ws.Cells.Item(131, 5).Value = $lessonName
$ws.Cells.Item(132, 5).Value = $lessonName
$ws.Cells.Item(133, 5).Value = $lessonName

# 2) "abordagem da aula" column (F), entered bottom-up (133 -> 130)
$ws.Cells.Item(133, 6).Value = "3:38 - fazer login no heroku CLI pelo terminaol cmd do windows -`nabrir CMD e digitar ""heroku login"" sem aspas`nvai pedir o login e senha da conta do heroku`n"
$ws.Cells.Item(132, 6).Value = "`n2:28 - instalação do software ""Heroku CLI"" no computador - programa que permite logar no heroku e efetuar operações remotas na aplicação"
$ws.Cells.Item(131, 6).Value = "1:30 - instalação do MySQL no heroku - adiciona uma instancia do MySQL ... necessário inserir um cartao de crédito`naba Overview->Configure Add-ons`npesquisar por ""MySQL"" e escolher o ClearDB MySQL`nEscolher um plano e clocar em ""Provide"""
$ws.Cells.Item(130, 6).Value = "0:50 - procedimento de criar app novo no heroku -`nCreatre new app`nnome do app (opcional)`nlocalidade (pais)"

# 3) Remaining columns - all reuse already-existing shared strings / numbers
$ws.Cells.Item(130, 2).Value = 4
$ws.Cells.Item(130, 3).Value = $sectionName
$ws.Cells.Item(130, 4).Value = 58
$ws.Cells.Item(130, 7).Value = "`n`n`n`n"

$ws.Cells.Item(131, 2).Value = 4
$ws.Cells.Item(131, 3).Value = $sectionName
$ws.Cells.Item(131, 4).Value = 58

$ws.Cells.Item(132, 2).Value = 4
$ws.Cells.Item(132, 3).Value = $sectionName
$ws.Cells.Item(132, 4).Value = 58

$ws.Cells.Item(133, 2).Value = 4
$ws.Cells.Item(133, 3).Value = $sectionName
$ws.Cells.Item(133, 4).Value = 58

# Copy formatting from the row above (129) which carries the styles used for
# this section (s=8/9/10 -> fonts/wrap/alignment), so the new rows/table rows
# pick up matching direct formatting instead of plain defaults.
for ($r = 130; $r -le 133; $r++) {
    $ws.Range("B129:G129").Copy() | Out-Null
    $ws.Range("B$r`:G$r").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

# Match row heights (content wraps inside the wide columns E/F/G, producing
# these auto heights in the original edit).
$ws.Rows.Item(130).RowHeight = 75
$ws.Rows.Item(131).RowHeight = 90
$ws.Rows.Item(132).RowHeight = 60
$ws.Rows.Item(133).RowHeight = 75

# Restore view/selection to match where the author ended up after the edit
$ws.Range("C127").Select()
$excel.ActiveWindow.ScrollRow = 127
$excel.ActiveWindow.ScrollColumn = 1
